$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted numeric-looking values stay as text (match original string cell type)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.208.56"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.270.95"
$ws.Range("E3").Value = "  -2.84%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.08"
$ws.Range("E5").Value = "  -3.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.29"
$ws.Range("E6").Value = "  -5.70%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.493"
$ws.Range("E8").Value = "  -3.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  -3.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.46"
$ws.Range("E10").Value = "  -3.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.23"
$ws.Range("E12").Value = "  -7.98%  "
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.66"
$ws.Range("E14").Value = "  -2.56%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.623.33"
$ws.Range("E15").Value = "  -3.07%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.54"
$ws.Range("E16").Value = "  -2.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.274.50"
$ws.Range("E17").Value = "  -2.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.782"
$ws.Range("E18").Value = "  -5.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.123.02"
$ws.Range("E19").Value = "  -1.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.71"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0890"
$ws.Range("E21").Value = "  -2.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.00"
$ws.Range("E22").Value = "  -2.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.51"
$ws.Range("E23").Value = "  -3.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "234.25"
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.97"
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.46"
$ws.Range("E27").Value = "  -4.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.94"
$ws.Range("E28").Value = "  -6.55%  "
$ws.Range("E29").Value = "  -1.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.62"
$ws.Range("E30").Value = "  +4.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.17"
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.10"
$ws.Range("E32").Value = "  -1.93%  "
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.90"
$ws.Range("E34").Value = "  -3.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.56"
$ws.Range("E35").Value = "  -1.34%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.53"
$ws.Range("E36").Value = "  -4.58%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.33"
$ws.Range("E37").Value = "  -4.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0688"
$ws.Range("E38").Value = "  -5.01%  "
$ws.Range("E39").Value = "  -3.97%  "
$ws.Range("E40").Value = "  -3.12%  "
$ws.Range("E41").Value = "  -3.26%  "
$ws.Range("E42").Value = "  -6.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.47"
$ws.Range("E43").Value = "  -4.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.962.67"
$ws.Range("E44").Value = "  -2.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0278"
$ws.Range("E45").Value = "  -2.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.43"
$ws.Range("E46").Value = "  -7.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.56"
$ws.Range("E47").Value = "  -6.62%  "
$ws.Range("E48").Value = "  -4.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.495.29"
$ws.Range("E49").Value = "  -2.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.13"
$ws.Range("E50").Value = "  -6.41%  "
$ws.Range("E51").Value = "  -3.60%  "
